$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(2,3,4,5,6,7,8,28,29,30,31,32,33,34,54,55,56,57,58,59,60,80,81,82,106,107,108,132,133,134)

foreach ($r in $rows) {
    $cell = $ws.Range("G$r")
    $val = $cell.Value2
    $parts = $val -split ',\s*'
    $n = $parts.Count
    $revParts = @()
    for ($i = $n - 1; $i -ge 0; $i--) {
        $revParts += $parts[$i]
    }
    $newVal = [string]::Join(', ', $revParts)
    $cell.Value = $newVal
}
